# Updates the crypto price/1h-volume columns (D, E) with freshly scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ D = <new price text, or $null if unchanged>; E = <new 1h-change text> }
$updates = @{
    2 = @{ D = '69.739.48'; E = '  -0.60%  ' }
    3 = @{ D = '3.833.19'; E = '  +2.53%  ' }
    4 = @{ D = '0.999'; E = '  -0.12%  ' }
    5 = @{ D = '611.34'; E = '  -1.94%  ' }
    6 = @{ D = '174.65'; E = '  -3.28%  ' }
    7 = @{ D = '3.829.44'; E = '  +2.51%  ' }
    8 = @{ D = $null; E = '  -0.01%  ' }
    9 = @{ D = '0.525'; E = '  -1.92%  ' }
    10 = @{ D = '0.166'; E = '  -1.50%  ' }
    11 = @{ D = '6.44'; E = '  +2.32%  ' }
    12 = @{ D = '0.478'; E = '  -1.95%  ' }
    13 = @{ D = '39.88'; E = '  -2.59%  ' }
    14 = @{ D = '0.0000253'; E = '  -2.63%  ' }
    15 = @{ D = '4.467.08'; E = '  +2.41%  ' }
    16 = @{ D = '3.829.52'; E = '  +2.50%  ' }
    17 = @{ D = '69.789.30'; E = '  -0.54%  ' }
    18 = @{ D = '7.44'; E = '  -2.21%  ' }
    19 = @{ D = $null; E = '  -3.44%  ' }
    20 = @{ D = '16.58'; E = '  -1.57%  ' }
    21 = @{ D = '504.91'; E = '  -0.28%  ' }
    22 = @{ D = '9.48'; E = '  +1.37%  ' }
    23 = @{ D = '0.734'; E = '  +1.04%  ' }
    24 = @{ D = '2.45'; E = '  -5.28%  ' }
    25 = @{ D = '85.72'; E = '  -1.16%  ' }
    26 = @{ D = $null; E = '  +3.51%  ' }
    27 = @{ D = '12.60'; E = '  -4.43%  ' }
    28 = @{ D = '10.39'; E = '  -9.31%  ' }
    29 = @{ D = $null; E = '  +0.18%  ' }
    30 = @{ D = '2.52'; E = '  +0.93%  ' }
    31 = @{ D = '2.98'; E = '  +1.04%  ' }
    32 = @{ D = '7.96'; E = '  -0.23%  ' }
    33 = @{ D = '32.08'; E = '  +2.71%  ' }
    34 = @{ D = $null; E = '  -2.52%  ' }
    35 = @{ D = '0.998'; E = '  -0.17%  ' }
    36 = @{ D = $null; E = '  -2.23%  ' }
    37 = @{ D = '6.07'; E = '  -2.26%  ' }
    38 = @{ D = '0.141'; E = '  +2.53%  ' }
    39 = @{ D = '487.67'; E = '  +13.66%  ' }
    40 = @{ D = '0.335'; E = '  -0.81%  ' }
    41 = @{ D = '2.05'; E = '  -3.15%  ' }
    42 = @{ D = '49.68'; E = '  -1.51%  ' }
    43 = @{ D = $null; E = '  +3.50%  ' }
    44 = @{ D = '43.24'; E = '  -5.54%  ' }
    45 = @{ D = '8.50'; E = '  -2.82%  ' }
    46 = @{ D = '2.920.13'; E = '  -2.92%  ' }
    47 = @{ D = $null; E = '  -1.66%  ' }
    48 = @{ D = '140.02'; E = '  +2.16%  ' }
    49 = @{ D = $null; E = '  +0.02%  ' }
    50 = @{ D = '26.77'; E = '  -2.70%  ' }
    51 = @{ D = '2.42'; E = '  -4.42%  ' }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($null -ne $vals.D) {
        $priceCell = $ws.Range("D$row")
        if ($vals.D -match '^[+-]?\d+(\.\d+)?$') {
            # Plain-numeric-looking text (e.g. "0.999") would otherwise be auto-converted
            # to a number on assignment; force Text so it is stored the same way the
            # scraper originally wrote it, then drop the temporary format stamp.
            $priceCell.NumberFormat = "@"
            $priceCell.Value = $vals.D
            $priceCell.ClearFormats()
        } else {
            $priceCell.Value = $vals.D
        }
    }
    if ($null -ne $vals.E) {
        $ws.Range("E$row").Value = $vals.E
    }
}
